# Apply updated NATMI values (Jag1-Notch2) per Dr Hou advice
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 10.27784
$ws.Cells.Item(2, 8).Value = 30.83352
$ws.Cells.Item(2, 9).Value = 0.230301226653591
$ws.Cells.Item(2, 10).Value = 0.230301226653591
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 27.53580066666666
$ws.Cells.Item(2, 14).Value = 82.60740199999999
$ws.Cells.Item(2, 15).Value = 0.2054887285464767
$ws.Cells.Item(2, 16).Value = 0.2054887285464768
$ws.Cells.Item(2, 17).Value = 283.0085535238933
$ws.Cells.Item(2, 18).Value = 2547.07698171504
$ws.Cells.Item(2, 19).Value = 0.04732430624774037
$ws.Cells.Item(2, 20).Value = 0.04732430624774037
# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 10.27784
$ws.Cells.Item(3, 8).Value = 30.83352
$ws.Cells.Item(3, 9).Value = 0.230301226653591
$ws.Cells.Item(3, 10).Value = 0.230301226653591
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 27.50472933333333
$ws.Cells.Item(3, 14).Value = 82.514188
$ws.Cells.Item(3, 15).Value = 0.2052568555438283
$ws.Cells.Item(3, 16).Value = 0.2052568555438283
$ws.Cells.Item(3, 17).Value = 282.6892073313066
$ws.Cells.Item(3, 18).Value = 2544.20286598176
$ws.Cells.Item(3, 19).Value = 0.04727090561080258
$ws.Cells.Item(3, 20).Value = 0.04727090561080258
# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 10.27784
$ws.Cells.Item(4, 8).Value = 30.83352
$ws.Cells.Item(4, 9).Value = 0.230301226653591
$ws.Cells.Item(4, 10).Value = 0.230301226653591
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 45.12975566666668
$ws.Cells.Item(4, 14).Value = 135.389267
$ws.Cells.Item(4, 15).Value = 0.3367854170582615
$ws.Cells.Item(4, 16).Value = 0.3367854170582616
$ws.Cells.Item(4, 17).Value = 463.8364079810934
$ws.Cells.Item(4, 18).Value = 4174.527671829841
$ws.Cells.Item(4, 19).Value = 0.07756209466755887
$ws.Cells.Item(4, 20).Value = 0.07756209466755887
# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 10.27784
$ws.Cells.Item(5, 8).Value = 30.83352
$ws.Cells.Item(5, 9).Value = 0.230301226653591
$ws.Cells.Item(5, 10).Value = 0.230301226653591
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 33.831228
$ws.Cells.Item(5, 14).Value = 101.493684
$ws.Cells.Item(5, 15).Value = 0.2524689988514334
$ws.Cells.Item(5, 16).Value = 0.2524689988514334
$ws.Cells.Item(5, 17).Value = 347.71194838752
$ws.Cells.Item(5, 18).Value = 3129.40753548768
$ws.Cells.Item(5, 19).Value = 0.05814392012748917
$ws.Cells.Item(5, 20).Value = 0.05814392012748916
# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 14.56812733333333
$ws.Cells.Item(6, 8).Value = 43.704382
$ws.Cells.Item(6, 9).Value = 0.3264360600001921
$ws.Cells.Item(6, 10).Value = 0.326436060000192
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 27.53580066666666
$ws.Cells.Item(6, 14).Value = 82.60740199999999
$ws.Cells.Item(6, 15).Value = 0.2054887285464767
$ws.Cells.Item(6, 16).Value = 0.2054887285464768
$ws.Cells.Item(6, 17).Value = 401.1450503372848
$ws.Cells.Item(6, 18).Value = 3610.305453035563
$ws.Cells.Item(6, 19).Value = 0.06707893092116086
$ws.Cells.Item(6, 20).Value = 0.06707893092116086
# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 14.56812733333333
$ws.Cells.Item(7, 8).Value = 43.704382
$ws.Cells.Item(7, 9).Value = 0.3264360600001921
$ws.Cells.Item(7, 10).Value = 0.326436060000192
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 27.50472933333333
$ws.Cells.Item(7, 14).Value = 82.514188
$ws.Cells.Item(7, 15).Value = 0.2052568555438283
$ws.Cells.Item(7, 16).Value = 0.2052568555438283
$ws.Cells.Item(7, 17).Value = 400.6923991968684
$ws.Cells.Item(7, 18).Value = 3606.231592771816
$ws.Cells.Item(7, 19).Value = 0.06700323921175588
$ws.Cells.Item(7, 20).Value = 0.06700323921175588
# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 14.56812733333333
$ws.Cells.Item(8, 8).Value = 43.704382
$ws.Cells.Item(8, 9).Value = 0.3264360600001921
$ws.Cells.Item(8, 10).Value = 0.326436060000192
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 45.12975566666668
$ws.Cells.Item(8, 14).Value = 135.389267
$ws.Cells.Item(8, 15).Value = 0.3367854170582615
$ws.Cells.Item(8, 16).Value = 0.3367854170582616
$ws.Cells.Item(8, 17).Value = 657.4560270742215
$ws.Cells.Item(8, 18).Value = 5917.104243667995
$ws.Cells.Item(8, 19).Value = 0.1099389046100204
$ws.Cells.Item(8, 20).Value = 0.1099389046100204
# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 14.56812733333333
$ws.Cells.Item(9, 8).Value = 43.704382
$ws.Cells.Item(9, 9).Value = 0.3264360600001921
$ws.Cells.Item(9, 10).Value = 0.326436060000192
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 33.831228
$ws.Cells.Item(9, 14).Value = 101.493684
$ws.Cells.Item(9, 15).Value = 0.2524689988514334
$ws.Cells.Item(9, 16).Value = 0.2524689988514334
$ws.Cells.Item(9, 17).Value = 492.857637347032
$ws.Cells.Item(9, 18).Value = 4435.718736123287
$ws.Cells.Item(9, 19).Value = 0.08241498525725494
$ws.Cells.Item(9, 20).Value = 0.08241498525725492
# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 2.650137
$ws.Cells.Item(10, 8).Value = 7.950411
$ws.Cells.Item(10, 9).Value = 0.05938308067649115
$ws.Cells.Item(10, 10).Value = 0.05938308067649114
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 27.53580066666666
$ws.Cells.Item(10, 14).Value = 82.60740199999999
$ws.Cells.Item(10, 15).Value = 0.2054887285464767
$ws.Cells.Item(10, 16).Value = 0.2054887285464768
$ws.Cells.Item(10, 17).Value = 72.97364417135799
$ws.Cells.Item(10, 18).Value = 656.762797542222
$ws.Cells.Item(10, 19).Value = 0.01220255374538502
$ws.Cells.Item(10, 20).Value = 0.01220255374538502
# Row 11
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 2.650137
$ws.Cells.Item(11, 8).Value = 7.950411
$ws.Cells.Item(11, 9).Value = 0.05938308067649115
$ws.Cells.Item(11, 10).Value = 0.05938308067649114
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 27.50472933333333
$ws.Cells.Item(11, 14).Value = 82.514188
$ws.Cells.Item(11, 15).Value = 0.2052568555438283
$ws.Cells.Item(11, 16).Value = 0.2052568555438283
$ws.Cells.Item(11, 17).Value = 72.89130088125199
$ws.Cells.Item(11, 18).Value = 656.021707931268
$ws.Cells.Item(11, 19).Value = 0.01218878441216204
$ws.Cells.Item(11, 20).Value = 0.01218878441216204
# Row 12
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 2.650137
$ws.Cells.Item(12, 8).Value = 7.950411
$ws.Cells.Item(12, 9).Value = 0.05938308067649115
$ws.Cells.Item(12, 10).Value = 0.05938308067649114
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 45.12975566666668
$ws.Cells.Item(12, 14).Value = 135.389267
$ws.Cells.Item(12, 15).Value = 0.3367854170582615
$ws.Cells.Item(12, 16).Value = 0.3367854170582616
$ws.Cells.Item(12, 17).Value = 119.600035293193
$ws.Cells.Item(12, 18).Value = 1076.400317638737
$ws.Cells.Item(12, 19).Value = 0.01999935559183646
$ws.Cells.Item(12, 20).Value = 0.01999935559183646
# Row 13
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 2.650137
$ws.Cells.Item(13, 8).Value = 7.950411
$ws.Cells.Item(13, 9).Value = 0.05938308067649115
$ws.Cells.Item(13, 10).Value = 0.05938308067649114
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 33.831228
$ws.Cells.Item(13, 14).Value = 101.493684
$ws.Cells.Item(13, 15).Value = 0.2524689988514334
$ws.Cells.Item(13, 16).Value = 0.2524689988514334
$ws.Cells.Item(13, 17).Value = 89.65738907823601
$ws.Cells.Item(13, 18).Value = 806.916501704124
$ws.Cells.Item(13, 19).Value = 0.01499238692710762
$ws.Cells.Item(13, 20).Value = 0.01499238692710762
# Row 14
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 17.13170833333333
$ws.Cells.Item(14, 8).Value = 51.395125
$ws.Cells.Item(14, 9).Value = 0.3838796326697257
$ws.Cells.Item(14, 10).Value = 0.3838796326697257
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 27.53580066666666
$ws.Cells.Item(14, 14).Value = 82.60740199999999
$ws.Cells.Item(14, 15).Value = 0.2054887285464767
$ws.Cells.Item(14, 16).Value = 0.2054887285464768
$ws.Cells.Item(14, 17).Value = 471.7353057461388
$ws.Cells.Item(14, 18).Value = 4245.61775171525
$ws.Cells.Item(14, 19).Value = 0.07888293763219048
$ws.Cells.Item(14, 20).Value = 0.07888293763219048
# Row 15
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 17.13170833333333
$ws.Cells.Item(15, 8).Value = 51.395125
$ws.Cells.Item(15, 9).Value = 0.3838796326697257
$ws.Cells.Item(15, 10).Value = 0.3838796326697257
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 27.50472933333333
$ws.Cells.Item(15, 14).Value = 82.514188
$ws.Cells.Item(15, 15).Value = 0.2052568555438283
$ws.Cells.Item(15, 16).Value = 0.2052568555438283
$ws.Cells.Item(15, 17).Value = 471.2030007259444
$ws.Cells.Item(15, 18).Value = 4240.8270065335
$ws.Cells.Item(15, 19).Value = 0.07879392630910775
$ws.Cells.Item(15, 20).Value = 0.07879392630910775
# Row 16
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 17.13170833333333
$ws.Cells.Item(16, 8).Value = 51.395125
$ws.Cells.Item(16, 9).Value = 0.3838796326697257
$ws.Cells.Item(16, 10).Value = 0.3838796326697257
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 45.12975566666668
$ws.Cells.Item(16, 14).Value = 135.389267
$ws.Cells.Item(16, 15).Value = 0.3367854170582615
$ws.Cells.Item(16, 16).Value = 0.3367854170582616
$ws.Cells.Item(16, 17).Value = 773.1498112359307
$ws.Cells.Item(16, 18).Value = 6958.348301123376
$ws.Cells.Item(16, 19).Value = 0.1292850621888458
$ws.Cells.Item(16, 20).Value = 0.1292850621888458
# Row 17
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 17.13170833333333
$ws.Cells.Item(17, 8).Value = 51.395125
$ws.Cells.Item(17, 9).Value = 0.3838796326697257
$ws.Cells.Item(17, 10).Value = 0.3838796326697257
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 33.831228
$ws.Cells.Item(17, 14).Value = 101.493684
$ws.Cells.Item(17, 15).Value = 0.2524689988514334
$ws.Cells.Item(17, 16).Value = 0.2524689988514334
$ws.Cells.Item(17, 17).Value = 579.5867306545
$ws.Cells.Item(17, 18).Value = 5216.2805758905
$ws.Cells.Item(17, 19).Value = 0.09691770653958165
$ws.Cells.Item(17, 20).Value = 0.09691770653958164
